$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-27 03:08:12"
$wsZh.Range("G5").Value = "2016-01-27 03:09:17"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-27 03:08:26"
$wsDe.Range("G5").Value = "2016-01-27 03:09:42"
